{"js": "// \"fixing the blank page code\": the document started with a standalone\n// FirstParagraph-styled paragraph (date + name/title line) followed by a\n// BodyText-styled paragraph with the actual summary text. That leading\n// paragraph is removed entirely, and the paragraph that used to follow it\n// (the summary) becomes the new first paragraph, taking on the\n// \"FirstParagraph\" style so the layout/spacing that used to start the page\n// is preserved without the extra (now blank-page-causing) paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// First paragraph = the date/name line to remove.\n// Second paragraph = the summary line that should become the new first\n// paragraph, inheriting the \"FirstParagraph\" style.\nconst firstPara = paragraphs.items[0];\nconst secondPara = paragraphs.items[1];\n\nsecondPara.style = \"FirstParagraph\";\nfirstPara.delete();\n\nawait context.sync();\n", "ps1": "# \"fixing the blank page code\": the document opened with a standalone\n# \"FirstParagraph\"-styled paragraph (the date + name/title line) followed by\n# a \"BodyText\"-styled paragraph holding the actual summary text. Remove that\n# leading paragraph entirely, and promote the paragraph that follows it (the\n# summary) to be the new first paragraph by giving it the \"FirstParagraph\"\n# style, so the document no longer starts with the extra paragraph that was\n# causing the blank page.\n\n$d = $word.ActiveDocument\n\n$firstPara = $d.Paragraphs(1)\n$secondPara = $d.Paragraphs(2)\n\n$secondPara.Style = \"FirstParagraph\"\n$firstPara.Range.Delete()\n"}
